{"js": "// Replace each addition/subtraction equation in the table with its updated value.\n// Old -> New pairs are applied in document order: each original equation text is located\n// with body.search (exact match, case sensitive) and its single hit is replaced in place.\nconst pairs = [\n  [\"79-41=38\", \"55+39=94\"],\n  [\"81+3=84\", \"51-24=27\"],\n  [\"8+66=74\", \"93-91=2\"],\n  [\"21+61=82\", \"57-16=41\"],\n  [\"85-3=82\", \"59+37=96\"],\n  [\"49+17=66\", \"97+0=97\"],\n  [\"75-31=44\", \"9+48=57\"],\n  [\"83-54=29\", \"84-33=51\"],\n  [\"51-50=1\", \"13+1=14\"],\n  [\"70-67=3\", \"28-11=17\"],\n  [\"99-75=24\", \"92-25=67\"],\n  [\"18-13=5\", \"6+63=69\"],\n  [\"16+72=88\", \"35+12=47\"],\n  [\"29+61=90\", \"97-69=28\"],\n  [\"13+72=85\", \"67-43=24\"],\n  [\"34+30=64\", \"30-9=21\"],\n  [\"28+64=92\", \"63-11=52\"],\n  [\"31+52=83\", \"13+70=83\"],\n  [\"29+64=93\", \"84-77=7\"],\n  [\"4+78=82\", \"17+41=58\"],\n  [\"23+40=63\", \"92-41=51\"],\n  [\"98-61=37\", \"22+29=51\"],\n  [\"22+57=79\", \"38+37=75\"],\n  [\"35+62=97\", \"62-39=23\"],\n  [\"53-32=21\", \"6+63=69\"],\n  [\"18+20=38\", \"88-78=10\"],\n  [\"88-1=87\", \"89-67=22\"],\n  [\"55-7=48\", \"66+11=77\"],\n  [\"75-11=64\", \"37-33=4\"],\n  [\"28+66=94\", \"21+33=54\"],\n  [\"6+55=61\", \"92-84=8\"],\n  [\"85+14=99\", \"82+16=98\"],\n  [\"99-53=46\", \"93-46=47\"],\n  [\"86-60=26\", \"87-71=16\"],\n  [\"82-18=64\", \"63-19=44\"],\n  [\"52-1=51\", \"42-24=18\"],\n  [\"65-1=64\", \"10+35=45\"],\n  [\"24+61=85\", \"69-53=16\"],\n  [\"93-32=61\", \"31-4=27\"],\n  [\"25+67=92\", \"73+17=90\"],\n  [\"66-22=44\", \"57+34=91\"],\n  [\"49+2=51\", \"49+4=53\"],\n  [\"16+5=21\", \"12+54=66\"],\n  [\"71-62=9\", \"22+48=70\"],\n  [\"87-40=47\", \"39+44=83\"],\n  [\"26+7=33\", \"29-28=1\"],\n  [\"18+9=27\", \"51+7=58\"],\n  [\"84-26=58\", \"58-44=14\"],\n  [\"81+6=87\", \"97-49=48\"],\n  [\"92-63=29\", \"48-16=32\"],\n  [\"15+76=91\", \"79-18=61\"],\n  [\"86+7=93\", \"41+30=71\"],\n  [\"1+71=72\", \"50+6=56\"],\n  [\"52-15=37\", \"30+8=38\"],\n  [\"53+44=97\", \"50+14=64\"],\n  [\"37+39=76\", \"40-10=30\"],\n  [\"6+33=39\", \"9+70=79\"],\n  [\"70-53=17\", \"87-35=52\"],\n  [\"52-42=10\", \"31+32=63\"],\n  [\"93-81=12\", \"36-35=1\"],\n  [\"27+71=98\", \"68-7=61\"],\n  [\"41+53=94\", \"76+10=86\"],\n  [\"13+81=94\", \"36+28=64\"],\n  [\"12+22=34\", \"60+23=83\"],\n  [\"41+44=85\", \"84-23=61\"],\n  [\"95-87=8\", \"44+31=75\"],\n  [\"50-38=12\", \"71-46=25\"],\n  [\"13+73=86\", \"26-0=26\"],\n  [\"28+44=72\", \"61-24=37\"],\n  [\"18+10=28\", \"95-24=71\"],\n  [\"3+28=31\", \"89-65=24\"],\n  [\"14+23=37\", \"84-71=13\"],\n  [\"26+36=62\", \"80-6=74\"],\n  [\"58+35=93\", \"49-36=13\"],\n  [\"3+67=70\", \"8+33=41\"],\n  [\"5+33=38\", \"5+81=86\"],\n  [\"24+48=72\", \"94-5=89\"],\n  [\"68-19=49\", \"52+38=90\"],\n  [\"6+48=54\", \"79+3=82\"],\n  [\"11+86=97\", \"32+57=89\"],\n  [\"38-24=14\", \"93-50=43\"],\n  [\"6+40=46\", \"36+51=87\"],\n  [\"3+10=13\", \"22+66=88\"],\n  [\"41-15=26\", \"44-31=13\"],\n  [\"86+9=95\", \"0+35=35\"],\n  [\"97-68=29\", \"59+26=85\"],\n  [\"6+68=74\", \"10+87=97\"],\n  [\"71-7=64\", \"2+22=24\"],\n  [\"21+64=85\", \"14+60=74\"],\n  [\"31+25=56\", \"15+6=21\"],\n  [\"65-11=54\", \"77-53=24\"],\n  [\"70+26=96\", \"26+27=53\"],\n  [\"73-9=64\", \"54-39=15\"],\n  [\"19+22=41\", \"35+25=60\"],\n  [\"56-3=53\", \"73-1=72\"],\n  [\"78-60=18\", \"85-83=2\"],\n  [\"99-84=15\", \"52-8=44\"],\n  [\"36+53=89\", \"86+3=89\"],\n  [\"67-50=17\", \"47+41=88\"],\n  [\"14+47=61\", \"58+41=99\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each addition/subtraction equation in the table with its updated value.\n# Old -> New pairs are applied in document order using Find/Replace on the document Range,\n# wrapping disabled and MatchWholeWord enabled so each original equation text is matched exactly once.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"79-41=38\", \"55+39=94\"),\n    @(\"81+3=84\", \"51-24=27\"),\n    @(\"8+66=74\", \"93-91=2\"),\n    @(\"21+61=82\", \"57-16=41\"),\n    @(\"85-3=82\", \"59+37=96\"),\n    @(\"49+17=66\", \"97+0=97\"),\n    @(\"75-31=44\", \"9+48=57\"),\n    @(\"83-54=29\", \"84-33=51\"),\n    @(\"51-50=1\", \"13+1=14\"),\n    @(\"70-67=3\", \"28-11=17\"),\n    @(\"99-75=24\", \"92-25=67\"),\n    @(\"18-13=5\", \"6+63=69\"),\n    @(\"16+72=88\", \"35+12=47\"),\n    @(\"29+61=90\", \"97-69=28\"),\n    @(\"13+72=85\", \"67-43=24\"),\n    @(\"34+30=64\", \"30-9=21\"),\n    @(\"28+64=92\", \"63-11=52\"),\n    @(\"31+52=83\", \"13+70=83\"),\n    @(\"29+64=93\", \"84-77=7\"),\n    @(\"4+78=82\", \"17+41=58\"),\n    @(\"23+40=63\", \"92-41=51\"),\n    @(\"98-61=37\", \"22+29=51\"),\n    @(\"22+57=79\", \"38+37=75\"),\n    @(\"35+62=97\", \"62-39=23\"),\n    @(\"53-32=21\", \"6+63=69\"),\n    @(\"18+20=38\", \"88-78=10\"),\n    @(\"88-1=87\", \"89-67=22\"),\n    @(\"55-7=48\", \"66+11=77\"),\n    @(\"75-11=64\", \"37-33=4\"),\n    @(\"28+66=94\", \"21+33=54\"),\n    @(\"6+55=61\", \"92-84=8\"),\n    @(\"85+14=99\", \"82+16=98\"),\n    @(\"99-53=46\", \"93-46=47\"),\n    @(\"86-60=26\", \"87-71=16\"),\n    @(\"82-18=64\", \"63-19=44\"),\n    @(\"52-1=51\", \"42-24=18\"),\n    @(\"65-1=64\", \"10+35=45\"),\n    @(\"24+61=85\", \"69-53=16\"),\n    @(\"93-32=61\", \"31-4=27\"),\n    @(\"25+67=92\", \"73+17=90\"),\n    @(\"66-22=44\", \"57+34=91\"),\n    @(\"49+2=51\", \"49+4=53\"),\n    @(\"16+5=21\", \"12+54=66\"),\n    @(\"71-62=9\", \"22+48=70\"),\n    @(\"87-40=47\", \"39+44=83\"),\n    @(\"26+7=33\", \"29-28=1\"),\n    @(\"18+9=27\", \"51+7=58\"),\n    @(\"84-26=58\", \"58-44=14\"),\n    @(\"81+6=87\", \"97-49=48\"),\n    @(\"92-63=29\", \"48-16=32\"),\n    @(\"15+76=91\", \"79-18=61\"),\n    @(\"86+7=93\", \"41+30=71\"),\n    @(\"1+71=72\", \"50+6=56\"),\n    @(\"52-15=37\", \"30+8=38\"),\n    @(\"53+44=97\", \"50+14=64\"),\n    @(\"37+39=76\", \"40-10=30\"),\n    @(\"6+33=39\", \"9+70=79\"),\n    @(\"70-53=17\", \"87-35=52\"),\n    @(\"52-42=10\", \"31+32=63\"),\n    @(\"93-81=12\", \"36-35=1\"),\n    @(\"27+71=98\", \"68-7=61\"),\n    @(\"41+53=94\", \"76+10=86\"),\n    @(\"13+81=94\", \"36+28=64\"),\n    @(\"12+22=34\", \"60+23=83\"),\n    @(\"41+44=85\", \"84-23=61\"),\n    @(\"95-87=8\", \"44+31=75\"),\n    @(\"50-38=12\", \"71-46=25\"),\n    @(\"13+73=86\", \"26-0=26\"),\n    @(\"28+44=72\", \"61-24=37\"),\n    @(\"18+10=28\", \"95-24=71\"),\n    @(\"3+28=31\", \"89-65=24\"),\n    @(\"14+23=37\", \"84-71=13\"),\n    @(\"26+36=62\", \"80-6=74\"),\n    @(\"58+35=93\", \"49-36=13\"),\n    @(\"3+67=70\", \"8+33=41\"),\n    @(\"5+33=38\", \"5+81=86\"),\n    @(\"24+48=72\", \"94-5=89\"),\n    @(\"68-19=49\", \"52+38=90\"),\n    @(\"6+48=54\", \"79+3=82\"),\n    @(\"11+86=97\", \"32+57=89\"),\n    @(\"38-24=14\", \"93-50=43\"),\n    @(\"6+40=46\", \"36+51=87\"),\n    @(\"3+10=13\", \"22+66=88\"),\n    @(\"41-15=26\", \"44-31=13\"),\n    @(\"86+9=95\", \"0+35=35\"),\n    @(\"97-68=29\", \"59+26=85\"),\n    @(\"6+68=74\", \"10+87=97\"),\n    @(\"71-7=64\", \"2+22=24\"),\n    @(\"21+64=85\", \"14+60=74\"),\n    @(\"31+25=56\", \"15+6=21\"),\n    @(\"65-11=54\", \"77-53=24\"),\n    @(\"70+26=96\", \"26+27=53\"),\n    @(\"73-9=64\", \"54-39=15\"),\n    @(\"19+22=41\", \"35+25=60\"),\n    @(\"56-3=53\", \"73-1=72\"),\n    @(\"78-60=18\", \"85-83=2\"),\n    @(\"99-84=15\", \"52-8=44\"),\n    @(\"36+53=89\", \"86+3=89\"),\n    @(\"67-50=17\", \"47+41=88\"),\n    @(\"14+47=61\", \"58+41=99\"),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n\nWrite-Output \"Replaced $($pairs.Count) equations\"\n"}
